$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37983.332
$ws.Range("J3").Value = 37983.332
$ws.Range("L3").Value = 37983.332
$ws.Range("N3").Value = -38211.332

$ws.Range("H34").Value = 10470.818
$ws.Range("I34").Value = 2354.2856
$ws.Range("J34").Value = 24674.75
$ws.Range("K34").Value = 2354.2856
$ws.Range("L34").Value = 24674.75
$ws.Range("M34").Value = -2151.2856
$ws.Range("N34").Value = -25080.75

$ws.Range("H36").Value = 10470.818
$ws.Range("I36").Value = 2354.2856
$ws.Range("J36").Value = 24674.75
$ws.Range("K36").Value = 2354.2856
$ws.Range("L36").Value = 24674.75
$ws.Range("M36").Value = -1639.2856
$ws.Range("N36").Value = -26104.75

$ws.Range("H80").Value = 10910.32
$ws.Range("I80").Value = 7519.1875
$ws.Range("J80").Value = 16939
$ws.Range("K80").Value = 22557.5625
$ws.Range("L80").Value = 50817
$ws.Range("M80").Value = -21559.5625
$ws.Range("N80").Value = -52813

$ws.Range("H83").Value = 10910.32
$ws.Range("I83").Value = 7519.1875
$ws.Range("J83").Value = 16939
$ws.Range("K83").Value = 67672.6875
$ws.Range("L83").Value = 152451
$ws.Range("M83").Value = -62680.6875
$ws.Range("N83").Value = -162435

$ws.Range("H102").Value = 37983.332
$ws.Range("J102").Value = 37983.332
$ws.Range("L102").Value = 37983.332
$ws.Range("N102").Value = -44473.332

$ws.Range("H137").Value = 1842952.9
$ws.Range("I137").Value = 2156461.2
$ws.Range("J137").Value = 1469154.4
$ws.Range("K137").Value = 6469383.600000001
$ws.Range("L137").Value = 4407463.199999999
$ws.Range("M137").Value = -6466833.600000001
$ws.Range("N137").Value = -4412563.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1547766.1
$ws.Range("I32").Value = 1547766.1
$ws.Range("K32").Value = 1547766.1
$ws.Range("M32").Value = -1547479.1

$ws.Range("H45").Value = 940.4
$ws.Range("I45").Value = 866.3333
$ws.Range("J45").Value = 1051.5
$ws.Range("K45").Value = 866.3333
$ws.Range("L45").Value = 1051.5
$ws.Range("M45").Value = -489.3333
$ws.Range("N45").Value = -1805.5

$ws.Range("H61").Value = 438440.66
$ws.Range("I61").Value = 314121.22
$ws.Range("J61").Value = 722599.4
$ws.Range("K61").Value = 314121.22
$ws.Range("L61").Value = 722599.4
$ws.Range("M61").Value = -313909.22
$ws.Range("N61").Value = -723023.4

$ws.Range("H110").Value = 1347.25
$ws.Range("I110").Value = 1312.9166
$ws.Range("J110").Value = 1398.75
$ws.Range("K110").Value = 1312.9166
$ws.Range("L110").Value = 1398.75
$ws.Range("M110").Value = 732.0834
$ws.Range("N110").Value = -5488.75

$ws.Range("H122").Value = 8888.5
$ws.Range("J122").Value = 7777
$ws.Range("L122").Value = 23331
$ws.Range("N122").Value = -28231

$ws.Range("H132").Value = 38512.355
$ws.Range("I132").Value = 78804.62
$ws.Range("J132").Value = 3592.4
$ws.Range("K132").Value = 236413.86
$ws.Range("L132").Value = 10777.2
$ws.Range("M132").Value = -233883.86
$ws.Range("N132").Value = -15837.2

$ws.Range("H136").Value = 438440.66
$ws.Range("I136").Value = 314121.22
$ws.Range("J136").Value = 722599.4
$ws.Range("K136").Value = 942363.6599999999
$ws.Range("L136").Value = 2167798.2
$ws.Range("M136").Value = -939813.6599999999
$ws.Range("N136").Value = -2172898.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 51499.5
$ws.Range("I4").Value = 10165.833
$ws.Range("J4").Value = 113500
$ws.Range("K4").Value = 10165.833
$ws.Range("L4").Value = 113500
$ws.Range("M4").Value = -10053.833
$ws.Range("N4").Value = -113724

$ws.Range("H28").Value = 29990
$ws.Range("J28").Value = 29990
$ws.Range("L28").Value = 29990
$ws.Range("N28").Value = -30480

$ws.Range("H31").Value = 1686776.6
$ws.Range("I31").Value = 980.5714
$ws.Range("J31").Value = 2759555.8
$ws.Range("K31").Value = 980.5714
$ws.Range("L31").Value = 2759555.8
$ws.Range("M31").Value = -685.5714
$ws.Range("N31").Value = -2760145.8

$ws.Range("H34").Value = 1686776.6
$ws.Range("I34").Value = 980.5714
$ws.Range("J34").Value = 2759555.8
$ws.Range("K34").Value = 980.5714
$ws.Range("L34").Value = 2759555.8
$ws.Range("M34").Value = -778.5714
$ws.Range("N34").Value = -2759959.8

$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10754
$ws.Range("M38").ClearContents()

$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10422
$ws.Range("M46").ClearContents()

$ws.Range("H106").Value = 27671
$ws.Range("J106").Value = 27671
$ws.Range("L106").Value = 27671
$ws.Range("N106").Value = -30195

$ws.Range("H122").Value = 1853.2106
$ws.Range("I122").Value = 2676.5
$ws.Range("J122").Value = 938.44446
$ws.Range("K122").Value = 8029.5
$ws.Range("L122").Value = 2815.33338
$ws.Range("M122").Value = -5579.5
$ws.Range("N122").Value = -7715.33338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 856.9032
$ws.Range("I68").Value = 506.6
$ws.Range("J68").Value = 1185.3125
$ws.Range("K68").Value = 1519.8
$ws.Range("L68").Value = 3555.9375
$ws.Range("M68").Value = -708.8000000000002
$ws.Range("N68").Value = -5177.9375

$ws.Range("H71").Value = 856.9032
$ws.Range("I71").Value = 506.6
$ws.Range("J71").Value = 1185.3125
$ws.Range("K71").Value = 4559.400000000001
$ws.Range("L71").Value = 10667.8125
$ws.Range("M71").Value = -503.4000000000005
$ws.Range("N71").Value = -18779.8125

$ws.Range("H107").Value = 1159.0322
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 1322.0834
$ws.Range("K107").Value = 1800
$ws.Range("L107").Value = 3966.2502
$ws.Range("M107").Value = 120
$ws.Range("N107").Value = -7806.2502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 33500
$ws.Range("J138").Value = 33500
$ws.Range("L138").Value = 33500
$ws.Range("N138").Value = -43780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 11142.223
$ws.Range("J104").Value = 11142.223
$ws.Range("L104").Value = 11142.223
$ws.Range("N104").Value = -18130.223

$ws.Range("H105").Value = 35999.5
$ws.Range("J105").Value = 35999.5
$ws.Range("L105").Value = 35999.5
$ws.Range("N105").Value = -42987.5

$ws.Range("H106").Value = 12494.75
$ws.Range("J106").Value = 12494.75
$ws.Range("L106").Value = 12494.75
$ws.Range("N106").Value = -15018.75

$ws.Range("H122").Value = 1490
$ws.Range("I122").Value = 1490
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4470
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("M122").Value = -2020

$ws.Range("H132").Value = 7582736.5
$ws.Range("I132").Value = 2560.074
$ws.Range("J132").Value = 19621840
$ws.Range("K132").Value = 7680.222
$ws.Range("L132").Value = 58865520
$ws.Range("M132").Value = -5150.222
$ws.Range("N132").Value = -58870580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 11597.777
$ws.Range("I43").Value = 10433.333
$ws.Range("K43").Value = 10433.333
$ws.Range("M43").Value = -10284.333

$ws.Range("H103").Value = 35500
$ws.Range("J103").Value = 35500
$ws.Range("L103").Value = 35500
$ws.Range("N103").Value = -37844
